# Insert a new data row at row 17 (pushes the existing rows 17-120 down to
# 18-121, growing the sheet's used range from A1:T120 to A1:T121) and fill it
# with a new weekly price observation for "Terminal Hortofrutícola Agro
# Chillán" / Mango.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..120 down by one row.
$ws.Rows("17:17").Insert()

# Populate the newly-inserted row 17 with the new record.
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 45035
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108002
$ws.Range("J17").Value = "Mango"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("Q17").Value = "$/bandeja 4 kilos"
$ws.Range("R17").Value = "Perú"
$ws.Range("S17").Value = 2000
$ws.Range("T17").Value = 4
